{"js": "// The document stores small inline-XML markup tags (e.g. \"<id>...</id>\")\n// as literal text split across several runs (one run for \"<id>\", one for\n// the value, one for \"</id>\") so each fragment can carry its own\n// formatting. For the \"p167r_1\" id, those three fragments should be\n// collapsed into a single run (keeping the opening tag's formatting)\n// whose text is the full \"<id>p167r_1</id>\" string.\nconst body = context.document.body;\n\n// This exact string is unique in the document: the sibling id\n// \"<id>fig_p167r_1</id>\" (and \"<id>fig_p167r_2</id>\") contain a\n// different value, so they are not matched and stay untouched.\nconst searchResults = body.search(\"<id>p167r_1</id>\", { matchCase: true, matchWildcards: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target text '<id>p167r_1</id>' not found in document.\");\n}\n\n// Replacing the whole found range with its own text merges the three\n// runs that made up the match into a single run, adopting the\n// formatting of the first run of the match (the \"<id>\" run), exactly\n// like the target edit.\nconst target = searchResults.items[0];\ntarget.insertText(\"<id>p167r_1</id>\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document stores small inline-XML markup tags (e.g. \"<id>...</id>\")\n# as literal text split across several runs (one run for the opening\n# \"<id>\" tag, one for the value, one for the closing \"</id>\" tag) so each\n# fragment can carry its own formatting. For the \"p167r_1\" id, those three\n# fragments should be collapsed into a single run (keeping the opening\n# tag's formatting) whose text is the full \"<id>p167r_1</id>\" string.\n\n# This exact literal string is unique in the document: the sibling ids\n# \"<id>fig_p167r_1</id>\" and \"<id>fig_p167r_2</id>\" contain a different\n# value, so they are not matched and stay untouched.\n$matchRange = $d.Content\n$find = $matchRange.Find\n$find.ClearFormatting()\n$find.Text = \"<id>p167r_1</id>\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$found = $find.Execute()\nif (-not $found) {\n  throw \"Target text '<id>p167r_1</id>' not found in document.\"\n}\n\n$matchStart = $matchRange.Start\n$matchEnd = $matchRange.End\n\n# Keep the formatting of the opening \"<id>\" run: only delete and\n# re-insert the remainder (\"p167r_1</id>\") right after it, which merges\n# everything into the first run instead of creating a new, unformatted\n# one.\n$openTag = \"<id>\"\n$splitPoint = $matchStart + $openTag.Length\n\n$headRange = $d.Range($matchStart, $splitPoint)\n$tailRange = $d.Range($splitPoint, $matchEnd)\n$tailText = $tailRange.Text\n$tailRange.Delete()\n$headRange.InsertAfter($tailText)\n"}
